$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.621.63"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "2.350.01"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.67"
$ws.Range("E5").Value = "  -3.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.27"
$ws.Range("E6").Value = "  +3.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  -2.06%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").Value = "  -5.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.50"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0926"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.62"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.00"
$ws.Range("E13").Value = "  -4.69%  "
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.00"
$ws.Range("E15").Value = "  -7.05%  "
$ws.Range("D16").Value = "2.703.25"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").Value = "2.356.01"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "42.578.99"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.70"
$ws.Range("E19").Value = "  -1.61%  "
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.94"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.71"
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.97"
$ws.Range("E23").Value = "  -6.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.32"
$ws.Range("E24").Value = "  -3.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.43"
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.43"
$ws.Range("E27").Value = "  -2.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.80"
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.94"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.21"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.06"
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.88"
$ws.Range("E34").Value = "  -8.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.126"
$ws.Range("E35").Value = "  +18.99%  "
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.66"
$ws.Range("E37").Value = "  -5.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0363"
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.93"
$ws.Range("E39").Value = "  -5.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.67"
$ws.Range("E40").Value = "  -4.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.240"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.48"
$ws.Range("E42").Value = "  -6.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.45"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.15"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "111.46"
$ws.Range("E46").Value = "  -8.41%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.22"
$ws.Range("E47").Value = "  -6.36%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.19"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.47"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.89"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("E51").Value = "  -2.09%  "
